# Automatische test-sync: 2025-08-06 20:21:50
# Appends a new log row to the "Logs" sheet and updates the "Dashboard"
# summary count for the matching category.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 16

$logs.Cells.Item($newRow, 1).Value = "Laat maar weten of er nieuws is"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Laat maar weten of er nieuws is"
$logs.Cells.Item($newRow, 4).Value = "Klantenservice / Opvolging"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-06 20:21:18"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional formatting ranges on the Logs sheet so they keep
# covering the newly added row (D/G/H/I/J 2:15 -> 2:16), preserving each
# rule's dxfId/priority/formula by re-pointing the sqref in place rather
# than deleting and re-adding the rules.
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "15")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "16")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard pivot-style summary: the "Klantenservice / Opvolging"
# category count goes from 1 to 2 now that a second row matches it.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(5, 2).Value = 2
